$d = $word.ActiveDocument
$r = $d.Range(49, 51)
$r.Text = "wAAAasasaw"
